$d = $word.ActiveDocument

# The paragraph currently holds the tag "<id>p123r_2</id>" split across
# three runs:
#   run 1 (Courier New, color 7f6000): "<id>"
#   run 2 (plain formatting):          "p123r_2"
#   run 3 (Courier New, color 7f6000): "</id>"
# The edit merges them into a single run (keeping run 1's Courier New
# "tag" formatting) whose text is "<id>p123r_2</id>".
#
# Strategy: delete runs 2+3's text ("p123r_2</id>") and re-append it to
# the end of run 1 ("<id>"). Extending run 1 in place (rather than
# building a brand new run) keeps its original run-level attributes.

# Delete "p123r_2</id>" (covers the 2nd and 3rd runs).
$rngTail = $d.Content
$rngTail.Find.ClearFormatting()
$foundTail = $rngTail.Find.Execute("p123r_2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundTail) {
    throw "Could not find 'p123r_2</id>' to remove."
}
$rngTail.Delete()

# Re-insert the removed text right after the remaining "<id>" run.
$rngHead = $d.Content
$rngHead.Find.ClearFormatting()
$foundHead = $rngHead.Find.Execute("<id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundHead) {
    throw "Could not find '<id>' run to extend."
}
$rngHead.InsertAfter("p123r_2</id>")
